# Set Runmode (column C) to "Y" for all test case rows in the "Test Cases" sheet,
# so the whole suite of authoring test cases runs.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$rows = @(2, 3, 4, 5, 6, 7, 8, 9, 14, 15, 20, 32)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = "Y"
}
